# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Recomputed "K" values (column G) for each row of match data, replacing the
# previously stored Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = @(4,0,1,1,1,1,1,1,3,1,2,1,0,0,0,1,2,2,1,4,4,0,1,0,1,2,0,2,1,2,1,1,2,1,2,1,2,1,0,4,3,2,0,0,2,0,0,3,0,1,1,0,3,3,4,1,1,2,1,2,2,0,2,3,3,2)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
